$wb = $excel.ActiveWorkbook

# ---- Regular_Section_A ----
$ws = $wb.Worksheets.Item("Regular_Section_A")
$ws.Range("D6").Value = "EC161 (Lab) [L105]"
$ws.Range("D7").Value = "EC161 (Lab) [L105]"
$ws.Range("E25").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D26").Value = "Mon 09:00-10:30 [C102], Wed 09:00-10:30 [C102]"
$ws.Range("E26").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D27").Value = "Mon 09:00-10:30 [C104], Wed 09:00-10:30 [C104]"
$ws.Range("E27").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 09:00-10:30 [C202], Wed 09:00-10:30 [C202]"
$ws.Range("E28").Value = "Fri 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 09:00-10:30 [C203], Wed 09:00-10:30 [C203]"
$ws.Range("E29").Value = "Fri 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 09:00-10:30 [C204], Wed 09:00-10:30 [C204]"
$ws.Range("E30").Value = "Fri 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 09:00-10:30 [C205], Wed 09:00-10:30 [C205]"
$ws.Range("E31").Value = "Fri 14:30-15:30 [C205]"
$ws.Range("D32").Value = "Mon 09:00-10:30 [C302], Wed 09:00-10:30 [C302]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C302]"
$ws.Range("D33").Value = "Mon 09:00-10:30 [C303], Wed 09:00-10:30 [C303]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C303]"
$ws.Range("D34").Value = "Mon 09:00-10:30 [C304], Wed 09:00-10:30 [C304]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C304]"

# ---- Regular_Section_B ----
$ws = $wb.Worksheets.Item("Regular_Section_B")
$ws.Range("C3").Value = "MA161 [C004]"
$ws.Range("B4").Value = "MA161 [C004]"
$ws.Range("C4").Value = "DS161 [C004]"
$ws.Range("D4").Value = "EC161 [C004]"
$ws.Range("B6").Value = "DS161 [C004]"
$ws.Range("C6").Value = "MA162 [C004]"
$ws.Range("D6").Value = "EC161 (Lab) [L206]"
$ws.Range("E6").Value = "CS161 (Lab) [L107]"
$ws.Range("D7").Value = "EC161 (Lab) [L206]"
$ws.Range("E7").Value = "CS161 (Lab) [L107]"
$ws.Range("B8").Value = "MA162 [C004]"
$ws.Range("C8").Value = "EC161 [C004]"
$ws.Range("E25").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D26").Value = "Mon 09:00-10:30 [C102], Wed 09:00-10:30 [C102]"
$ws.Range("E26").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D27").Value = "Mon 09:00-10:30 [C104], Wed 09:00-10:30 [C104]"
$ws.Range("E27").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 09:00-10:30 [C202], Wed 09:00-10:30 [C202]"
$ws.Range("E28").Value = "Fri 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 09:00-10:30 [C203], Wed 09:00-10:30 [C203]"
$ws.Range("E29").Value = "Fri 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 09:00-10:30 [C204], Wed 09:00-10:30 [C204]"
$ws.Range("E30").Value = "Fri 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 09:00-10:30 [C205], Wed 09:00-10:30 [C205]"
$ws.Range("E31").Value = "Fri 14:30-15:30 [C205]"
$ws.Range("D32").Value = "Mon 09:00-10:30 [C302], Wed 09:00-10:30 [C302]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C302]"
$ws.Range("D33").Value = "Mon 09:00-10:30 [C303], Wed 09:00-10:30 [C303]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C303]"
$ws.Range("D34").Value = "Mon 09:00-10:30 [C304], Wed 09:00-10:30 [C304]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C304]"

# ---- PreMid_Section_A ----
$ws = $wb.Worksheets.Item("PreMid_Section_A")
$ws.Range("D4").Value = "MA161 [C004]"
$ws.Range("C6").Value = "MA161 [C004]"
$ws.Range("C8").Value = "DS161 [C004]"
$ws.Range("E25").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D26").Value = "Mon 09:00-10:30 [C102], Wed 09:00-10:30 [C102]"
$ws.Range("E26").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D27").Value = "Mon 09:00-10:30 [C104], Wed 09:00-10:30 [C104]"
$ws.Range("E27").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 09:00-10:30 [C202], Wed 09:00-10:30 [C202]"
$ws.Range("E28").Value = "Fri 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 09:00-10:30 [C203], Wed 09:00-10:30 [C203]"
$ws.Range("E29").Value = "Fri 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 09:00-10:30 [C204], Wed 09:00-10:30 [C204]"
$ws.Range("E30").Value = "Fri 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 09:00-10:30 [C205], Wed 09:00-10:30 [C205]"
$ws.Range("E31").Value = "Fri 14:30-15:30 [C205]"
$ws.Range("D32").Value = "Mon 09:00-10:30 [C302], Wed 09:00-10:30 [C302]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C302]"
$ws.Range("D33").Value = "Mon 09:00-10:30 [C303], Wed 09:00-10:30 [C303]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C303]"
$ws.Range("D34").Value = "Mon 09:00-10:30 [C304], Wed 09:00-10:30 [C304]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C304]"

# ---- PreMid_Section_B ----
$ws = $wb.Worksheets.Item("PreMid_Section_B")
$ws.Range("B4").Value = "CS161 [C101]"
$ws.Range("D4").Value = "MA161 [C004]"
$ws.Range("C6").Value = "MA161 [C004]"
$ws.Range("D6").Value = "DS161 [C004]"
$ws.Range("B8").Value = "HS161 [C101]"
$ws.Range("C8").Value = "DS161 [C004]"
$ws.Range("E25").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D26").Value = "Mon 09:00-10:30 [C102], Wed 09:00-10:30 [C102]"
$ws.Range("E26").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D27").Value = "Mon 09:00-10:30 [C104], Wed 09:00-10:30 [C104]"
$ws.Range("E27").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 09:00-10:30 [C202], Wed 09:00-10:30 [C202]"
$ws.Range("E28").Value = "Fri 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 09:00-10:30 [C203], Wed 09:00-10:30 [C203]"
$ws.Range("E29").Value = "Fri 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 09:00-10:30 [C204], Wed 09:00-10:30 [C204]"
$ws.Range("E30").Value = "Fri 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 09:00-10:30 [C205], Wed 09:00-10:30 [C205]"
$ws.Range("E31").Value = "Fri 14:30-15:30 [C205]"
$ws.Range("D32").Value = "Mon 09:00-10:30 [C302], Wed 09:00-10:30 [C302]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C302]"
$ws.Range("D33").Value = "Mon 09:00-10:30 [C303], Wed 09:00-10:30 [C303]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C303]"
$ws.Range("D34").Value = "Mon 09:00-10:30 [C304], Wed 09:00-10:30 [C304]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C304]"

# ---- PostMid_Section_A ----
$ws = $wb.Worksheets.Item("PostMid_Section_A")
$ws.Range("B4").Value = "CS161 [C102]"
$ws.Range("D4").Value = "MA162 [C004]"
$ws.Range("C6").Value = "MA162 [C004]"
$ws.Range("D6").Value = "EC161 [C004]"
$ws.Range("B8").Value = "HS161 [C102]"
$ws.Range("C8").Value = "EC161 [C004]"
$ws.Range("D8").Value = "EC161 (Lab) [L105]"
$ws.Range("D9").Value = "EC161 (Lab) [L105]"
$ws.Range("E25").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D26").Value = "Mon 09:00-10:30 [C102], Wed 09:00-10:30 [C102]"
$ws.Range("E26").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D27").Value = "Mon 09:00-10:30 [C104], Wed 09:00-10:30 [C104]"
$ws.Range("E27").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 09:00-10:30 [C202], Wed 09:00-10:30 [C202]"
$ws.Range("E28").Value = "Fri 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 09:00-10:30 [C203], Wed 09:00-10:30 [C203]"
$ws.Range("E29").Value = "Fri 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 09:00-10:30 [C204], Wed 09:00-10:30 [C204]"
$ws.Range("E30").Value = "Fri 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 09:00-10:30 [C205], Wed 09:00-10:30 [C205]"
$ws.Range("E31").Value = "Fri 14:30-15:30 [C205]"
$ws.Range("D32").Value = "Mon 09:00-10:30 [C302], Wed 09:00-10:30 [C302]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C302]"
$ws.Range("D33").Value = "Mon 09:00-10:30 [C303], Wed 09:00-10:30 [C303]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C303]"
$ws.Range("D34").Value = "Mon 09:00-10:30 [C304], Wed 09:00-10:30 [C304]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C304]"

# ---- PostMid_Section_B ----
$ws = $wb.Worksheets.Item("PostMid_Section_B")
$ws.Range("B4").Value = "CS161 [C104]"
$ws.Range("D4").Value = "MA162 [C004]"
$ws.Range("C6").Value = "MA162 [C004]"
$ws.Range("D6").Value = "EC161 [C004]"
$ws.Range("B8").Value = "HS161 [C104]"
$ws.Range("C8").Value = "EC161 [C004]"
$ws.Range("D8").Value = "EC161 (Lab) [L206]"
$ws.Range("D9").Value = "EC161 (Lab) [L206]"
$ws.Range("E25").Value = "Fri 14:30-15:30 [C101]"
$ws.Range("D26").Value = "Mon 09:00-10:30 [C102], Wed 09:00-10:30 [C102]"
$ws.Range("E26").Value = "Fri 14:30-15:30 [C102]"
$ws.Range("D27").Value = "Mon 09:00-10:30 [C104], Wed 09:00-10:30 [C104]"
$ws.Range("E27").Value = "Fri 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 09:00-10:30 [C202], Wed 09:00-10:30 [C202]"
$ws.Range("E28").Value = "Fri 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 09:00-10:30 [C203], Wed 09:00-10:30 [C203]"
$ws.Range("E29").Value = "Fri 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 09:00-10:30 [C204], Wed 09:00-10:30 [C204]"
$ws.Range("E30").Value = "Fri 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 09:00-10:30 [C205], Wed 09:00-10:30 [C205]"
$ws.Range("E31").Value = "Fri 14:30-15:30 [C205]"
$ws.Range("D32").Value = "Mon 09:00-10:30 [C302], Wed 09:00-10:30 [C302]"
$ws.Range("E32").Value = "Fri 14:30-15:30 [C302]"
$ws.Range("D33").Value = "Mon 09:00-10:30 [C303], Wed 09:00-10:30 [C303]"
$ws.Range("E33").Value = "Fri 14:30-15:30 [C303]"
$ws.Range("D34").Value = "Mon 09:00-10:30 [C304], Wed 09:00-10:30 [C304]"
$ws.Range("E34").Value = "Fri 14:30-15:30 [C304]"
